$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '56.575.99'
$ws.Range("E2").Value = '  -4.30%  '

# Row 3
$ws.Range("D3").Value = '2.375.42'
$ws.Range("E3").Value = '  -5.68%  '

# Row 4
$ws.Range("E4").Value = '  +0.23%  '

# Row 5
$ws.Range("D5").Value = '''510.01'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -4.74%  '

# Row 6
$ws.Range("D6").Value = '''130.17'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.19%  '

# Row 7
$ws.Range("D7").Value = '''0.998'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.19%  '

# Row 8
$ws.Range("D8").Value = '''0.555'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.55%  '

# Row 9
$ws.Range("D9").Value = '2.399.59'
$ws.Range("E9").Value = '  -4.69%  '

# Row 10
$ws.Range("D10").Value = '''0.0964'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.15%  '

# Row 11
$ws.Range("E11").Value = '  -1.79%  '

# Row 12
$ws.Range("B12").Value = 'Cardano'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D12").Value = '''0.321'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -3.45%  '

# Row 13
$ws.Range("B13").Value = 'Toncoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D13").Value = '''4.69'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -10.32%  '

# Row 14
$ws.Range("D14").Value = '2.800.57'
$ws.Range("E14").Value = '  -5.40%  '

# Row 15
$ws.Range("D15").Value = '56.458.18'
$ws.Range("E15").Value = '  -4.35%  '

# Row 16
$ws.Range("D16").Value = '''21.69'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.49%  '

# Row 18
$ws.Range("D18").Value = '2.378.05'
$ws.Range("E18").Value = '  -5.39%  '

# Row 19
$ws.Range("D19").Value = '''10.25'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -4.01%  '

# Row 20
$ws.Range("D20").Value = '''4.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.45%  '

# Row 21
$ws.Range("D21").Value = '''312.31'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.91%  '

# Row 22
$ws.Range("D22").Value = '''6.27'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.00%  '

# Row 23
$ws.Range("D23").Value = '''0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.26%  '

# Row 24
$ws.Range("D24").Value = '''65.38'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.10%  '

# Row 25
$ws.Range("D25").Value = '''0.998'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.15%  '

# Row 26
$ws.Range("B26").Value = 'WrappedeETH'
$ws.Range("C26").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D26").Value = '2.484.28'
$ws.Range("E26").Value = '  -5.49%  '

# Row 27
$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").Value = '''0.379'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -7.95%  '

# Row 28
$ws.Range("E28").Value = '  -4.72%  '

# Row 29
$ws.Range("D29").Value = '''7.26'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.60%  '

# Row 30
$ws.Range("D30").Value = '''174.92'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.52%  '

# Row 32
$ws.Range("D32").Value = '0.0₃0718'
$ws.Range("E32").Value = '  -5.71%  '

# Row 33
$ws.Range("E33").Value = '  -2.64%  '

# Row 34
$ws.Range("E34").Value = '  -5.97%  '

# Row 35
$ws.Range("D35").Value = '''0.998'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '

# Row 36
$ws.Range("E36").Value = '  -0.39%  '

# Row 37
$ws.Range("E37").Value = '  -2.54%  '

# Row 38
$ws.Range("E38").Value = '  -3.77%  '

# Row 39
$ws.Range("D39").Value = '''3.74'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.21%  '

# Row 40
$ws.Range("D40").Value = '''35.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.26%  '

# Row 41
$ws.Range("D41").Value = '''1.43'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -5.49%  '

# Row 42
$ws.Range("D42").Value = '''0.790'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -3.38%  '

# Row 43
$ws.Range("D43").Value = '''129.95'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.64%  '

# Row 44
$ws.Range("E44").Value = '  -4.16%  '

# Row 45
$ws.Range("E45").Value = '  -3.87%  '

# Row 46
$ws.Range("E46").Value = '  -3.29%  '

# Row 47
$ws.Range("D47").Value = '''256.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -7.67%  '

# Row 48
$ws.Range("E48").Value = '  -3.75%  '

# Row 49
$ws.Range("D49").Value = '''0.0490'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.55%  '

# Row 50
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").Value = '''0.0208'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -5.55%  '

# Row 51
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '''16.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.36%  '
